$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") stores a date serial that is bumped by exactly one
# day (2023-09-12 -> 2023-09-13, i.e. serial 45181 -> 45182) for every data
# row in the sheet (rows 2 through 391). Row 1 is the header and is left
# untouched.
$ws.Range("C2:C391").Value = 45182
